$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (06-January-2025): the day was marked "At Work" (C7=1) but should
# instead be marked "Sick Leave" (E7=1), with a matching remark in H7.

# Move the highlighted/bordered cell format: C7 should become a plain cell
# (copy the formatting already used by the plain "Public Holiday" cell D7),
# and E7 should become the highlighted cell (copy the formatting already
# used by the highlighted "At Work" cell C2).
$ws.Range("D7").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# Swap the values
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 1

# Add the remark
$ws.Range("H7").Value = "Sick Leave"

# Update the Total row (row 33) to reflect the change: one less "At Work"
# day and one more "Sick Leave" day.
$ws.Range("C33").Value = 22
$ws.Range("E33").Value = 1
